$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.730.13'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '3.126.61'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''532.73'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('D6').Value = '''138.38'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.123.66'
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('E9').Value = '  +6.52%  '
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('E12').Value = '  +4.94%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '3.656.26'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').Value = '''26.02'
$ws.Range('E15').Value = '  +3.10%  '
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '57.819.88'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '3.123.85'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '''6.08'
$ws.Range('E19').Value = '  +2.76%  '
$ws.Range('D20').Value = '''12.72'
$ws.Range('E20').Value = '  +2.95%  '
$ws.Range('D21').Value = '''8.08'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('D22').Value = '''368.34'
$ws.Range('E22').Value = '  +6.94%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('D25').Value = '''69.17'
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').Value = '''0.506'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('D27').Value = '''0.168'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '0.0₃0866'
$ws.Range('E29').Value = '  -2.50%  '
$ws.Range('D30').Value = '''7.30'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').Value = '''6.07'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').Value = '''21.42'
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('D34').Value = '''5.15'
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').Value = '''159.48'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('D38').Value = '''1.30'
$ws.Range('E38').Value = '  +5.67%  '
$ws.Range('D39').Value = '''25.46'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').Value = '''1.68'
$ws.Range('E40').Value = '  +4.41%  '
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').Value = '2.547.21'
$ws.Range('E42').Value = '  +7.27%  '
$ws.Range('D43').Value = '''4.08'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '''0.699'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').Value = '''37.79'
$ws.Range('E45').Value = '  +3.50%  '
$ws.Range('D48').Value = '''0.976'
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('D49').Value = '''6.11'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').Value = '''19.73'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').Value = '''0.740'
$ws.Range('E51').Value = '  -2.18%  '

# Row 46/47 swap (VeChain <-> FirstDigitalUSD)
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '''0.0269'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '''1.00'
$ws.Range('E47').Value = '  -0.14%  '
